# Guia M1 sin soluciones
# -----------------------------------------------------------------------
# 1) Hoja1: clear the "answer" formatting from the FILTER helper block so
#    it reads as a blank worksheet, move the selection, and
# 2) add a second sheet "Solución" (after Hoja1) that holds the fully
#    worked example (values + SUMIF/COUNTIF/FILTER formulas + scratch
#    columns), and make it the active sheet/tab.
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

# ---- currency (Moneda) number format string, reused several times -----
$moneda = '_-"$"* #,##0.00_-;\-"$"* #,##0.00_-;_-"$"* "-"??_-;_-@_-'

# =========================================================================
# 1) Hoja1 clean-up
# =========================================================================

# Strip the shading/border formatting off the worked "USANDO FILTRAR"
# answer block (F8:I10 and F13:I15) -- values/formulas are untouched.
$src.Range("F8:I8").ClearFormats()
$src.Range("G9:I9").ClearFormats()
$src.Range("F13:I13").ClearFormats()
$src.Range("G14:I14").ClearFormats()

# F9 and F14 were empty placeholder cells that only carried formatting --
# remove them entirely.
$src.Range("F9").ClearContents()
$src.Range("F9").ClearFormats()
$src.Range("F14").ClearContents()
$src.Range("F14").ClearFormats()

# Move the selection/active cell (was J11) and drop this sheet as the
# tab that's flagged "selected" -- Solución becomes the active tab below.
$src.Range("D12").Select()

# =========================================================================
# 2) Add the "Solución" worksheet right after Hoja1
# =========================================================================

$sol = $wb.Worksheets.Add($null, $src)
$sol.Name = "Solución"

# ---- base table (A1:D7), plain values, no fills/borders ---------------
$sol.Range("A1").Value = "Cliente"
$sol.Range("B1").Value = "Género"
$sol.Range("C1").Value = "Crédito"
$sol.Range("D1").Value = "Monto"

$sol.Range("A2:D7").Value2 = $src.Range("A2:D7").Value2
$sol.Range("D2:D7").NumberFormat = $moneda

# ---- SUMIF-based summary (F2:G3) ---------------------------------------
$sol.Range("F2").Value = "Hombre"
$sol.Range("G2").Formula = '=SUMIF(B2:B7,"Hombre",D2:D7)'
$sol.Range("G2").NumberFormat = $moneda

$sol.Range("F3").Value = "Mujer"
$sol.Range("G3").Formula = '=SUMIF(B2:B7,"Mujer",D2:D7)'
$sol.Range("G3").NumberFormat = $moneda

# ---- FILTER-based breakdown (F5:J7) ------------------------------------
$sol.Range("F5").Value = "Hombre"
$sol.Range("G5:J7").FormulaArray = '=_xlfn._xlws.FILTER(A2:D7,B2:B7="Hombre")'

$sol.Range("F6").Value = "Mujer"

# ---- scratch/demo columns (D10:E15, G10:G12, I10) ----------------------
$sol.Range("D10:D15").FormulaArray = "=B2:B7"
$sol.Range("E10:E15").FormulaArray = '=B2:B7="Hombre"'
$sol.Range("E10:E15").NumberFormat = $moneda

$sol.Range("G10:G12").FormulaArray = '=_xlfn._xlws.FILTER(D2:D7,B2:B7="Mujer")'

$sol.Range("I10").FormulaArray = '=SUM(_xlfn._xlws.FILTER(D2:D7,B2:B7="Hombre"))'
$sol.Range("I10").NumberFormat = $moneda

# ---- cosmetic column widths (bestFit in the source file) --------------
$sol.Columns.Item(4).ColumnWidth = 10.6667
$sol.Columns.Item(5).ColumnWidth = 12
$sol.Columns.Item(7).ColumnWidth = 10.6667
$sol.Columns.Item(9).ColumnWidth = 10.6667

# ---- activate Solución, zoom + selection -------------------------------
$sol.Activate()
$excel.ActiveWindow.Zoom = 230
$sol.Range("B15").Select()
